# "upgrade left table until javakheti"
# C. Poti - Number of stillbirths
#   - rename the worksheet to match the place name
#   - mark the "Rural" row's 2010-2020 figures as confidential/unavailable ("...")
#   - drop the blank separator row above the footnote

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was generically named "1" - give it the proper title.
$ws.Name = "C.Poti"

# Row 7 ("Rural") previously showed 0 for 2010-2020 (cols B-L); these values
# are not actually known/published, so replace them with the same
# confidential/unavailable marker already used elsewhere in the row (M:O).
$ws.Range("B7:L7").Value = "..."

# Remove the empty row 8 that separated the data table from the footnote,
# shifting the footnote row (formerly row 9) up to row 8.
$ws.Rows("8").Delete()
